# Figure_17.4-5.pptx -- "Fix figures in Sections 10 through 17"
#
# Repositions four callout text boxes on slide 1 (Form Filler, Form
# Receiver, Retrieve Clarifications [ITI-37], Form Source) and refreshes
# the cached "Date Updates Automatically" footer field (datetimeFigureOut)
# on the slide master and every slide layout from 8/12/20 to 8/26/20.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Move/resize the four callout text boxes on the (single) slide.
#    Coordinates below are expressed in points (PowerPoint COM's native
#    unit for Shape.Left/Top/Width/Height -- 1 pt = 12700 EMU) and were
#    chosen so that they round-trip to the exact target EMU values.
# ---------------------------------------------------------------------

# "Form Filler" text box -> only the vertical position changes
# (436563,495300) -> (436563,751777) EMU
$shpFiller = $s.Shapes.Item(1)
$shpFiller.Left = 34.37504007007874
$shpFiller.Top  = 59.195039770078736

# "Form Receiver" text box -> both x & y shift slightly
# (3267075,752475) -> (3345132,763626) EMU
$shpReceiver = $s.Shapes.Item(5)
$shpReceiver.Left = 263.396225072441
$shpReceiver.Top  = 60.128031496062995

# "Retrieve Clarifications [ITI-37]" text box -> moves down/right
# (1139825,984250) -> (1173278,1151515) EMU
$shpClarify = $s.Shapes.Item(7)
$shpClarify.Left = 92.38409808818898
$shpClarify.Top  = 90.67047504094488

# "Form Source" text box -> moves and widens
# off (2790825,495300) -> (2822695,439661) EMU
# ext (914400,295275)  -> (1257059,295275) EMU
$shpSource = $s.Shapes.Item(14)
$shpSource.Left   = 222.25945291889764
$shpSource.Top    = 34.618976377952755
$shpSource.Width  = 98.98102362204725
$shpSource.Height = 23.25

# ---------------------------------------------------------------------
# 2. Refresh the cached auto-date footer text (8/12/20 -> 8/26/20) on the
#    slide master and on every slide layout.
# ---------------------------------------------------------------------

function Update-CachedDate {
    param($shapeRange)

    for ($i = 1; $i -le $shapeRange.Count; $i++) {
        $shp = $shapeRange.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "8/12/20") {
                $shp.TextFrame.TextRange.Text = "8/26/20"
            }
        }
    }
}

Update-CachedDate $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-CachedDate $p.SlideMaster.CustomLayouts.Item($li).Shapes
}
